# Add a "WEEK 3" column to the team-contributions table.
#
# The worksheet has a single Excel Table ("Table") spanning B2:D6
# (MEMBERS, WEEK 1, WEEK 2). We extend it one column to the right with a
# new "WEEK 3" header and fill in the same contribution numbers as the
# other weeks (33/33/33/0), matching the pattern already used for
# WEEK 1 / WEEK 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing table and grow it by one column (B2:D6 -> B2:E6).
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()

# The cell just to the right of the table (E4) carries the thin "table
# glow" border left over from the original layout; when the table grows,
# that stray formatting shifts one column further right (to F4). Preserve
# it by copying E4's current format onto F4 before we overwrite E4 with
# real data.
$ws.Range("E4").Copy()
$ws.Range("F4").PasteSpecial(-4122)

# Header + data for the new "WEEK 3" column.
$ws.Range("E2").Value = "WEEK 3"
$ws.Range("E2").Font.Size = 18
$ws.Range("E3").Value = 33
$ws.Range("E4").Value = 33
$ws.Range("E5").Value = 33
$ws.Range("E6").Value = 0

# E4 should end up with plain formatting (no border) now that it holds
# real data - the border moved to F4 above.
$ws.Range("E4").Borders.LineStyle = -4142

# Match the column width used for the new column.
$ws.Columns.Item(5).ColumnWidth = 15.5

# Leave the selection where it ended up after entering the data.
[void]$ws.Range("F7").Select()
